$wb = $excel.ActiveWorkbook

# --- ShipmentInformation sheet: update pickup/dropoff reference codes ---
$wsShip = $wb.Worksheets.Item("ShipmentInformation")
$wsShip.Range("C2").Value = "PickUp174"
$wsShip.Range("K2").Value = "DropOff53"

# --- Input sheet: update claim/tracking numbers for the claim row ---
$wsInput = $wb.Worksheets.Item("Input")

# T3 holds a purely numeric looking claim id that must stay text (shared string),
# so it's entered with a leading quote to force text storage, then the cell's
# visual formatting (thin top/bottom border + solid fill) is restored.
$wsInput.Range("T3").Value = "'58285547"
$wsInput.Range("T3").Borders.Item(8).LineStyle = 1
$wsInput.Range("T3").Borders.Item(9).LineStyle = 1
$wsInput.Range("T3").Interior.ColorIndex = 2

$wsInput.Range("W3").Value = "FCT913438719413846016"
$wsInput.Range("X3").Value = "FCTEST1003998"

# B3, U3 and Y3 did not change value but were re-touched as part of the same
# claim-row refresh, so re-apply their (unchanged) formatting.
$wsInput.Range("B3").Borders.Item(8).LineStyle = 1
$wsInput.Range("B3").Borders.Item(9).LineStyle = 1
$wsInput.Range("B3").Interior.ColorIndex = 2

$wsInput.Range("U3").Borders.Item(8).LineStyle = 1
$wsInput.Range("U3").Borders.Item(9).LineStyle = 1
$wsInput.Range("U3").Interior.ColorIndex = 2

$wsInput.Range("Y3").Borders.Item(8).LineStyle = 1
$wsInput.Range("Y3").Borders.Item(9).LineStyle = 1
$wsInput.Range("Y3").Interior.ColorIndex = 2

# --- ClaimDetail sheet: re-touched as part of the same refresh ---
$wsClaim = $wb.Worksheets.Item("ClaimDetail")
$wsClaim.Range("C3").Borders.Item(8).LineStyle = 1
$wsClaim.Range("C3").Borders.Item(9).LineStyle = 1
$wsClaim.Range("C3").Interior.ColorIndex = 2
